# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$oldText = $wsHoja1.Range("A1").Value()
$newText = $oldText
$newText = $newText.Replace("1000 Bs = 2.08 = 7529.54 pesos", "1000 Bs = 2.01 = 7229.21 pesos")
$newText = $newText.Replace("7529.54 pesos = 2.08 = 949.04 Bs", "7229.21 pesos = 2.0 = 916.56 Bs")
$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update the N10/O10/N12/O12 rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 497.98
$wsTasas.Range("O10").Value = 3600
$wsTasas.Range("N12").Value = 3620.9
$wsTasas.Range("O12").Value = 459.08
